$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15: "echo >>" entry
$ws.Range("A15").Value = "echo >>"
$ws.Range("B15").Value = "attach content to a file"
$ws.Range("C15").Value = 'echo "Another line" >> another-file.txt'

# Row 4: add a third column example for "nano" (typed after the row-15 block)
$ws.Range("C4").Value = "nano new-file.txt"

# New row 16: "rm" entry (remove files)
$ws.Range("A16").Value = "rm"
$ws.Range("C16").Value = "rm new-file.txt"
$ws.Range("B16").Value = "remove files"

# New row 17: "rm" entry (remove directories)
$ws.Range("A17").Value = "rm"
$ws.Range("B17").Value = "directories"
$ws.Range("C17").Value = "rm -rf new-folder"

# Update the active selection to mirror the author's final cursor position
$ws.Range("B22").Select() | Out-Null
